$wb = $excel.ActiveWorkbook

# --- Sheet "Layer0" (first sheet) ---
$ws1 = $wb.Worksheets.Item("Layer0")

$ws1.Range("B2").Value = 6.473123644325607
$ws1.Range("C2").Value = -4.607381193671313
$ws1.Range("D2").Value = -1.141886234607804
$ws1.Range("E2").Value = 3.606202155191703

$ws1.Range("B3").Value = 3.311868186720423
$ws1.Range("C3").Value = -2.068866931281668
$ws1.Range("D3").Value = -1.253255737341218
$ws1.Range("E3").Value = 1.278720260749817

$ws1.Range("B4").Value = -10.43802135758207
$ws1.Range("C4").Value = -3.538907799284799
$ws1.Range("D4").Value = 1.460551027431167
$ws1.Range("E4").Value = -1.925301796362968

$ws1.Range("B5").Value = 0.8198212490174183
$ws1.Range("C5").Value = 2.678582464998535
$ws1.Range("D5").Value = -0.8722260191830133
$ws1.Range("E5").Value = -3.03262737927183

$ws1.Range("B6").Value = 4.549684261512178
$ws1.Range("C6").Value = 4.443183458040439
$ws1.Range("D6").Value = 3.215193690146705
$ws1.Range("E6").Value = -1.609457440631691

$ws1.Range("B7").Value = 0.6853265773749967
$ws1.Range("C7").Value = 1.118973986793525
$ws1.Range("D7").Value = 0.6874921761912348
$ws1.Range("E7").Value = 7.749149725933131

# --- Sheet "Layer1" (second sheet) ---
$ws2 = $wb.Worksheets.Item("Layer1")

$ws2.Range("B2").Value = -18.99172214959324
$ws2.Range("C2").Value = -4.122332012459756
$ws2.Range("D2").Value = -0.7321956849236436
$ws2.Range("E2").Value = -11.03407452550561
$ws2.Range("F2").Value = 1.180184025604254

$ws2.Range("B3").Value = 1.629410162324189
$ws2.Range("C3").Value = 7.225536175766389
$ws2.Range("D3").Value = -9.590779546743681
$ws2.Range("E3").Value = 8.560315622877614
$ws2.Range("F3").Value = -3.690752232856454

$ws2.Range("B4").Value = 7.237419701901674
$ws2.Range("C4").Value = 7.835217450405891
$ws2.Range("D4").Value = -4.879155978740124
$ws2.Range("E4").Value = -9.563612158599508
$ws2.Range("F4").Value = -8.200205498644948

$ws2.Range("B5").Value = 8.944480098602378
$ws2.Range("C5").Value = -11.25671029181228
$ws2.Range("D5").Value = -8.71683329159104
$ws2.Range("E5").Value = -2.258608497314561
$ws2.Range("F5").Value = 11.08050166900075

$ws2.Range("B6").Value = 8.990569027750178
$ws2.Range("C6").Value = -9.528284903960067
$ws2.Range("D6").Value = 5.179210161542684
$ws2.Range("E6").Value = 7.285550227792794
$ws2.Range("F6").Value = -7.886814118333672
